# Auto-generated Excel COM-interop script to apply F-column (registration/visitor count)
# updates across sheets "展览", "演出", and "全部类型".

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (30 updates)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 17
$ws.Range("F3").Value = 26
$ws.Range("F5").Value = 5100
$ws.Range("F6").Value = 5100
$ws.Range("F7").Value = 103
$ws.Range("F9").Value = 512
$ws.Range("F11").Value = 1150
$ws.Range("F13").Value = 4956
$ws.Range("F15").Value = 60
$ws.Range("F16").Value = 75
$ws.Range("F18").Value = 221
$ws.Range("F20").Value = 242
$ws.Range("F21").Value = 3761
$ws.Range("F23").Value = 37
$ws.Range("F24").Value = 3645
$ws.Range("F25").Value = 171
$ws.Range("F28").Value = 209
$ws.Range("F30").Value = 200
$ws.Range("F35").Value = 136
$ws.Range("F36").Value = 6448
$ws.Range("F37").Value = 1025
$ws.Range("F38").Value = 484
$ws.Range("F39").Value = 96
$ws.Range("F42").Value = 1321
$ws.Range("F44").Value = 641
$ws.Range("F46").Value = 2218
$ws.Range("F47").Value = 311
$ws.Range("F48").Value = 88
$ws.Range("F49").Value = 762
$ws.Range("F50").Value = 907

# Sheet: 演出 (2 updates)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 4
$ws.Range("F23").Value = 801

# Sheet: 全部类型 (24 updates)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 17
$ws.Range("F6").Value = 4
$ws.Range("F8").Value = 26
$ws.Range("F10").Value = 5100
$ws.Range("F11").Value = 5100
$ws.Range("F12").Value = 103
$ws.Range("F17").Value = 4956
$ws.Range("F19").Value = 60
$ws.Range("F20").Value = 75
$ws.Range("F23").Value = 242
$ws.Range("F24").Value = 3645
$ws.Range("F25").Value = 171
$ws.Range("F27").Value = 209
$ws.Range("F29").Value = 200
$ws.Range("F34").Value = 136
$ws.Range("F36").Value = 6448
$ws.Range("F37").Value = 1025
$ws.Range("F38").Value = 96
$ws.Range("F40").Value = 1321
$ws.Range("F42").Value = 641
$ws.Range("F44").Value = 2218
$ws.Range("F45").Value = 311
$ws.Range("F47").Value = 88
$ws.Range("F48").Value = 762

